# Update the crypto price/volume table in-place to match the
# "Updated cryptos list" GitHub Actions refresh.
#
# The Price (D) and Volume(1h) (E) columns are stored as literal text
# (t="inlineStr") even though most Price values look like plain numbers
# ("1.008", "19.40", ...). A plain `Range.Value = "1.008"` assignment lets
# Excel's COM layer auto-coerce that into a real number, which would change
# the cell's stored type. To keep these as text (matching the workbook's
# original authoring), force NumberFormat to Text ("@") before the write,
# then restore the cell style to Normal afterwards so no stray formatting
# is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextCell "D2" "26.593.18"
Set-TextCell "E2" "  -2.21%  "

# Row 3 - Ethereum
Set-TextCell "D3" "1.815.30"
Set-TextCell "E3" "  -2.01%  "

# Row 4 - TetherUSD
Set-TextCell "D4" "1.007"
Set-TextCell "E4" "  +0.58%  "

# Row 5 - was USDC, now BNB (rows 5/6 swapped rank order)
Set-TextCell "B5" "BNB"
Set-TextCell "C5" "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
Set-TextCell "D5" "308.56"
Set-TextCell "E5" "  -1.51%  "

# Row 6 - was BNB, now USDC
Set-TextCell "B6" "USDC"
Set-TextCell "C6" "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextCell "D6" "1.006"
Set-TextCell "E6" "  +0.49%  "

# Row 7 - XRP
Set-TextCell "D7" "0.4570"
Set-TextCell "E7" "  -1.55%  "

# Row 8 - Cardano
Set-TextCell "D8" "0.3670"
Set-TextCell "E8" "  -1.26%  "

# Row 9 - Dogecoin
Set-TextCell "D9" "0.07150"
Set-TextCell "E9" "  -1.80%  "

# Row 10 - Polygon
Set-TextCell "D10" "0.8807"
Set-TextCell "E10" "  -0.70%  "

# Row 11 - TRON
Set-TextCell "D11" "0.07789"
Set-TextCell "E11" "  -0.65%  "

# Row 12 - Solana
Set-TextCell "D12" "19.43"
Set-TextCell "E12" "  -2.97%  "

# Row 13 - WrappedEther
Set-TextCell "D13" "1.777.22"
Set-TextCell "E13" "  -6.95%  "

# Row 14 - Polkadot
Set-TextCell "D14" "5.298"
Set-TextCell "E14" "  -1.50%  "

# Row 15 - Chainlink (price unchanged, volume only)
Set-TextCell "E15" "  -2.10%  "

# Row 16 - Litecoin
Set-TextCell "D16" "86.33"
Set-TextCell "E16" "  -4.94%  "

# Row 17 - BinanceUSD (price unchanged, volume only)
Set-TextCell "E17" "  +0.67%  "

# Row 18 - ShibaInu
Set-TextCell "D18" "0.000008605"
Set-TextCell "E18" "  -3.50%  "

# Row 19 - Dai
Set-TextCell "D19" "1.007"
Set-TextCell "E19" "  +0.63%  "

# Row 20 - WrappedBTC
Set-TextCell "D20" "26.640.85"
Set-TextCell "E20" "  -2.15%  "

# Row 21 - Avalanche
Set-TextCell "D21" "14.30"
Set-TextCell "E21" "  -2.93%  "

# Row 22 - Uniswap
Set-TextCell "D22" "5.009"
Set-TextCell "E22" "  -1.16%  "

# Row 23 - Cosmos
Set-TextCell "D23" "10.47"
Set-TextCell "E23" "  -0.28%  "

# Row 24 - Toncoin (price unchanged, volume only)
Set-TextCell "E24" "  +1.51%  "

# Row 25 - Monero
Set-TextCell "D25" "150.86"
Set-TextCell "E25" "  -0.58%  "

# Row 26 - EthereumClassic
Set-TextCell "D26" "18.02"
Set-TextCell "E26" "  -1.96%  "

# Row 27 - LidoDAOToken
Set-TextCell "D27" "2.076"
Set-TextCell "E27" "  +1.89%  "

# Row 28 - BitcoinCash
Set-TextCell "D28" "112.90"
Set-TextCell "E28" "  -2.39%  "

# Row 29 - InternetComputer(DFINITY)
Set-TextCell "D29" "4.869"
Set-TextCell "E29" "  -3.66%  "

# Row 30 - Stellar
Set-TextCell "D30" "0.08687"
Set-TextCell "E30" "  -1.35%  "

# Row 31 - HuobiToken
Set-TextCell "D31" "3.048"
Set-TextCell "E31" "  -3.08%  "

# Row 32 - ImmutableX
Set-TextCell "D32" "0.7361"
Set-TextCell "E32" "  -3.90%  "

# Row 33 - Filecoin (volume unchanged, price only)
Set-TextCell "D33" "4.484"

# Row 34 - ARBITRUM
Set-TextCell "D34" "1.121"
Set-TextCell "E34" "  -3.87%  "

# Row 35 - was Frax, now RenderToken (rows 35/36 swapped rank order)
Set-TextCell "B35" "RenderToken"
Set-TextCell "C35" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D35" "2.590"
Set-TextCell "E35" "  -4.49%  "

# Row 36 - was RenderToken, now Frax
Set-TextCell "B36" "Frax"
Set-TextCell "C36" "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextCell "D36" "1.003"
Set-TextCell "E36" "  +0.32%  "

# Row 37 - TrustWalletToken
Set-TextCell "D37" "1.081"
Set-TextCell "E37" "  -2.98%  "

# Row 38 - VeChain: unchanged, no edits

# Row 39 - Hedera
Set-TextCell "D39" "0.05123"
Set-TextCell "E39" "  -1.46%  "

# Row 40 - MXToken
Set-TextCell "D40" "2.896"
Set-TextCell "E40" "  -1.38%  "

# Row 41 - FraxShare
Set-TextCell "D41" "6.983"
Set-TextCell "E41" "  -0.51%  "

# Row 42 - TheSandbox (volume unchanged, price only)
Set-TextCell "D42" "0.5019"

# Row 43 - Algorand (price unchanged, volume only)
Set-TextCell "E43" "  -3.83%  "

# Row 44 - Aptos
Set-TextCell "D44" "8.177"
Set-TextCell "E44" "  -2.93%  "

# Row 45 - PaxDollar
Set-TextCell "D45" "1.007"
Set-TextCell "E45" "  +0.64%  "

# Row 46 - Decentraland
Set-TextCell "D46" "0.4631"
Set-TextCell "E46" "  -3.32%  "

# Row 47 - EnergySwap
Set-TextCell "D47" "9.966"
Set-TextCell "E47" "  -3.76%  "

# Row 48 - Quant
Set-TextCell "D48" "101.10"
Set-TextCell "E48" "  -1.62%  "

# Row 49 - NEARProtocol (volume unchanged, price only)
Set-TextCell "D49" "1.596"

# Row 50 - Cronos
Set-TextCell "D50" "0.06027"
Set-TextCell "E50" "  -2.88%  "

# Row 51 - Aave
Set-TextCell "D51" "64.18"
Set-TextCell "E51" "  -2.04%  "
